$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of sample data to append to the bottom of the table (rows 10-14)
$rows = @(
    @("NOT_VALID", "Record no REC ID", "Record without REC ID", "Brad (14583316)", "Rejected", "Apple", "Orange", "Operational", "Record without REC ID procedure.", "Record without REC ID procedure.", "Record without REC ID procedure.", 3),
    @("NOT_VALID", "Another Record no REC ID", "Another record that doesn't have REC ID", "Brad (14583316)", "Rejected", "Apple", "Orange", "Operational", "Another record without REC ID procedure.", "Another record without REC ID procedure.", "Another record without REC ID procedure.", 4),
    @("REC-1111", "A retired record", "Some retired record description", "Brad (14583316)", "Retired", "Apple", "Orange", "Security", "Retired record procedure.", "Retired record procedure.", "Retired record procedure.", 3),
    @("REC-0000", "An orphaned record", "Orphaned record for sure", "Jackson (03582313)", "Active", "Apple", "Orange", "Security", "Orphaned record procedure.", "Orphaned record procedure.", "Orphaned record procedure.", 4),
    @("REC-0001", "An orphaned record", "Orphaned record for sure", "Jackson (03582313)", "Active", "Apple", "Orange", "Security", "Orphaned record procedure.", "Orphaned record procedure.", "Orphaned record procedure.", 4)
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}

# Widen column A slightly to fit the longer "Record ID" values now present
$ws.Columns.Item(1).ColumnWidth = 11.65

# Update the active selection to reflect where the editor left off
$ws.Range("A15:C16").Select()
